# Adding Account Exist Test cases
$wb = $excel.ActiveWorkbook

# --- Fix AddNewCustomerTest run-mode flag (last row): N -> Y ---
$addNewCustomerTest = $wb.Worksheets.Item("AddNewCustomerTest")
$addNewCustomerTest.Range("E5").Value = "Y"

# --- Insert a new "AccountExistTest" sheet right before "test_suite" ---
$testSuite = $wb.Worksheets.Item("test_suite")
$accountExistTest = $wb.Worksheets.Add($testSuite)
$accountExistTest.Name = "AccountExistTest"

$accountExistTest.Range("A1").Value = "fullName"
$accountExistTest.Range("B1").Value = "noAccountText"
$accountExistTest.Range("A2").Value = "Rounak Agarwal"
$accountExistTest.Range("B2").Value = "Please open an account with us."
$accountExistTest.Range("A3").Value = "Sapnish Singh"
$accountExistTest.Range("B3").Value = "Please open an account with us."
$accountExistTest.Range("A4").Value = "Kuntal Chakraborty"

$accountExistTest.Columns.Item(1).AutoFit() | Out-Null
$accountExistTest.Columns.Item(2).AutoFit() | Out-Null

$accountExistTest.Range("B8").Select() | Out-Null

# --- Restore the first sheet as the active/selected tab ---
$addNewCustomerTest.Select()
$addNewCustomerTest.Range("D8").Select() | Out-Null
